$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('PayNowCC')
$ws.Range('B2').Value = 'Wed Nov 05 23:26:07 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 23:27:02 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:27:53 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:28:46 IST 2025'
$ws.Range('B6').Value = 'Wed Nov 05 23:29:44 IST 2025'
$ws.Range('B7').Value = 'Wed Nov 05 23:30:38 IST 2025'
$ws.Range('B8').Value = 'Wed Nov 05 23:31:37 IST 2025'
$ws.Range('B9').Value = 'Wed Nov 05 23:32:29 IST 2025'
$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsCredit')
$ws.Range('B2').Value = 'Wed Nov 05 23:59:55 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:01:01 IST 2025'
$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsCredit')
$ws.Range('B2').Value = 'Thu Nov 06 00:26:58 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:27:57 IST 2025'
$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsPC')
$ws.Range('B2').Value = 'Thu Nov 06 00:03:56 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:05:15 IST 2025'
$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsPC')
$ws.Range('B2').Value = 'Thu Nov 06 00:30:48 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:31:58 IST 2025'
$ws = $wb.Worksheets.Item('PayNowCorpSCF')
$ws.Range('B2').Value = 'Wed Nov 05 22:47:36 IST 2025'
$ws.Range('A3').Value = 'Fail'
$ws.Range('B3').Value = 'Wed Nov 05 22:48:52 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 22:51:12 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 22:52:21 IST 2025'
$ws = $wb.Worksheets.Item('PayNowCreditSCF')
$ws.Range('A2').Value = 'Fail'
$ws.Range('B2').Value = 'Wed Nov 05 22:59:20 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 23:01:39 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:02:55 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:04:17 IST 2025'
$ws = $wb.Worksheets.Item('PayNowPersonalCheckSCF')
$ws.Range('B2').Value = 'Wed Nov 05 23:10:23 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 23:11:44 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:12:58 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:14:17 IST 2025'
$ws = $wb.Worksheets.Item('PayNowPersonalSavingsSCF')
$ws.Range('B2').Value = 'Wed Nov 05 23:15:35 IST 2025'
$ws.Range('A3').Value = 'Fail'
$ws.Range('B3').Value = 'Wed Nov 05 23:16:51 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:19:05 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:20:17 IST 2025'
$ws = $wb.Worksheets.Item('PayNowCreditDCF')
$ws.Range('B2').Value = 'Wed Nov 05 22:54:32 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 22:55:43 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 22:56:52 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 22:58:04 IST 2025'
$ws = $wb.Worksheets.Item('PayNowCorpDCF')
$ws.Range('B2').Value = 'Wed Nov 05 22:42:44 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 22:43:54 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 22:45:06 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 22:46:25 IST 2025'
$ws = $wb.Worksheets.Item('PayNowPC')
$ws.Range('B2').Value = 'Thu Nov 06 21:35:42 IST 2025'
$ws = $wb.Worksheets.Item('PayNowPersonalCheckDCF')
$ws.Range('B2').Value = 'Wed Nov 05 23:21:31 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 23:22:42 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:23:47 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:24:59 IST 2025'
$ws = $wb.Worksheets.Item('MaxAmountErrorCC')
$ws.Range('B2').Value = 'Thu Nov 06 01:18:40 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:19:39 IST 2025'
$ws = $wb.Worksheets.Item('MaxAmountErrorCorp')
$ws.Range('B2').Value = 'Thu Nov 06 01:20:32 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:21:26 IST 2025'
$ws = $wb.Worksheets.Item('MaxAmountErrorPC')
$ws.Range('B2').Value = 'Thu Nov 06 01:22:25 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:23:25 IST 2025'
$ws = $wb.Worksheets.Item('MaxAmountErrorPS')
$ws.Range('B2').Value = 'Thu Nov 06 01:24:25 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:25:24 IST 2025'
$ws = $wb.Worksheets.Item('MinAmountErrorPC')
$ws.Range('B2').Value = 'Thu Nov 06 01:30:09 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:31:12 IST 2025'
$ws = $wb.Worksheets.Item('MinAmountErrorCC')
$ws.Range('B2').Value = 'Thu Nov 06 01:26:18 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:27:13 IST 2025'
$ws = $wb.Worksheets.Item('MinAmountErrorCorp')
$ws.Range('B2').Value = 'Thu Nov 06 01:28:06 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:29:13 IST 2025'
$ws = $wb.Worksheets.Item('MinAmountErrorPS')
$ws.Range('B2').Value = 'Thu Nov 06 01:32:02 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:32:54 IST 2025'
$ws = $wb.Worksheets.Item('PayNowPS')
$ws.Range('B2').Value = 'Thu Nov 06 21:43:14 IST 2025'
$ws.Range('A3').Value = 'Fail'
$ws.Range('B3').Value = 'Wed Nov 05 22:37:58 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 22:39:19 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 22:40:09 IST 2025'
$ws.Range('B6').Value = 'Wed Nov 05 22:41:01 IST 2025'
$ws.Range('B7').Value = 'Wed Nov 05 22:41:51 IST 2025'
$ws = $wb.Worksheets.Item('OverAndUnderPayCredit')
$ws.Range('A2').Value = 'Fail'
$ws.Range('B2').Value = 'Wed Nov 05 23:38:08 IST 2025'
$ws.Range('A3').Value = 'Fail'
$ws.Range('B3').Value = 'Wed Nov 05 23:39:34 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:41:01 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:41:57 IST 2025'
$ws = $wb.Worksheets.Item('OverAndUnderPayPC')
$ws.Range('B2').Value = 'Wed Nov 05 23:42:49 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 23:43:44 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:44:41 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:45:38 IST 2025'
$ws = $wb.Worksheets.Item('OverAndUnderPayPS')
$ws.Range('B2').Value = 'Wed Nov 05 23:46:37 IST 2025'
$ws.Range('B3').Value = 'Wed Nov 05 23:47:28 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:48:24 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:49:16 IST 2025'
$ws = $wb.Worksheets.Item('OverAndUnderPayCorp')
$ws.Range('A2').Value = 'Fail'
$ws.Range('B2').Value = 'Wed Nov 05 23:33:20 IST 2025'
$ws.Range('A3').Value = 'Fail'
$ws.Range('B3').Value = 'Wed Nov 05 23:34:45 IST 2025'
$ws.Range('B4').Value = 'Wed Nov 05 23:36:08 IST 2025'
$ws.Range('B5').Value = 'Wed Nov 05 23:37:06 IST 2025'
$ws = $wb.Worksheets.Item('NoModifyAmountCorp')
$ws.Range('A2').Value = 'Fail'
$ws.Range('B2').Value = 'Thu Nov 06 01:48:32 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:50:08 IST 2025'
$ws = $wb.Worksheets.Item('NoModifyAmountPC')
$ws.Range('B2').Value = 'Thu Nov 06 01:51:08 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:52:12 IST 2025'
$ws = $wb.Worksheets.Item('NoModifyAmountPS')
$ws.Range('B2').Value = 'Thu Nov 06 01:53:07 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:54:00 IST 2025'
$ws = $wb.Worksheets.Item('NoModifyAmountCC')
$ws.Range('B2').Value = 'Thu Nov 06 01:46:48 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:47:43 IST 2025'
$ws = $wb.Worksheets.Item('NoOverPayErrorCC')
$ws.Range('B2').Value = 'Thu Nov 06 01:35:30 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:36:23 IST 2025'
$ws = $wb.Worksheets.Item('NoOverPayErrorPC')
$ws.Range('B2').Value = 'Thu Nov 06 01:39:18 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:40:18 IST 2025'
$ws = $wb.Worksheets.Item('PayNowCorp')
$ws.Range('B2').Value = 'Wed Nov 05 22:53:33 IST 2025'
$ws = $wb.Worksheets.Item('NoOverPayErrorCorp')
$ws.Range('B2').Value = 'Thu Nov 06 01:37:22 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:38:20 IST 2025'
$ws = $wb.Worksheets.Item('NoOverPayErrorPS')
$ws.Range('B2').Value = 'Thu Nov 06 01:41:12 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:42:02 IST 2025'
$ws = $wb.Worksheets.Item('NoUnderPayErrorPS')
$ws.Range('B2').Value = 'Thu Nov 06 01:45:50 IST 2025'
$ws = $wb.Worksheets.Item('NoUnderPayErrorPC')
$ws.Range('B2').Value = 'Thu Nov 06 01:44:51 IST 2025'
$ws = $wb.Worksheets.Item('NoUnderPayErrorCC')
$ws.Range('B2').Value = 'Thu Nov 06 01:42:53 IST 2025'
$ws = $wb.Worksheets.Item('NoUnderPayErrorCorp')
$ws.Range('B2').Value = 'Thu Nov 06 01:44:00 IST 2025'
$ws = $wb.Worksheets.Item('CardExpiredErrorCC')
$ws.Range('A2').Value = 'Pass'
$ws.Range('B2').Value = 'Thu Nov 06 01:14:58 IST 2025'
$ws.Range('A3').Value = 'Pass'
$ws.Range('B3').Value = 'Thu Nov 06 01:15:50 IST 2025'
$ws = $wb.Worksheets.Item('CardNotAcceptedErrorCC')
$ws.Range('B2').Value = 'Thu Nov 06 01:16:45 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:17:45 IST 2025'
$ws = $wb.Worksheets.Item('MRFCorp')
$ws.Range('B2').Value = 'Thu Nov 06 01:33:46 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 01:34:41 IST 2025'
$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsCorp')
$ws.Range('B2').Value = 'Thu Nov 06 00:02:00 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:03:00 IST 2025'
$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsCorp')
$ws.Range('B2').Value = 'Thu Nov 06 00:28:53 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:29:50 IST 2025'
$ws = $wb.Worksheets.Item('VerifyConfirmPageLabelsPS')
$ws.Range('B2').Value = 'Thu Nov 06 00:06:24 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:07:30 IST 2025'
$ws = $wb.Worksheets.Item('VerifyReceiptPageLabelsPS')
$ws.Range('B2').Value = 'Thu Nov 06 00:33:06 IST 2025'
$ws.Range('B3').Value = 'Thu Nov 06 00:34:15 IST 2025'
